# Maj diapo revue 1
# The deck's footer date placeholder (an automatically-updating
# "datetimeFigureOut" field cached as 25/01/2018) is refreshed to
# 29/01/2018 everywhere it is defined: the slide master, every custom
# (slide) layout, and the notes master.

$p = $ppt.ActivePresentation

$oldDate = "25/01/2018"
$newDate = "29/01/2018"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq $true -and $shp.TextFrame.HasText -eq $true) {
            $isDate = $false
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDate = $true
                }
            }
            if ($isDate -eq $true) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# 1) Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2) Every custom layout hanging off the slide master.
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# 3) The notes master.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
